# Add duplicate content to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right after row 2 (pushes the old row 3 down to row 6).
$ws.Range("A3:A5").EntireRow.Insert()

# Fill the newly inserted rows 3-5 with a copy of row 2's data.
$ws.Range("A2:C2").Copy($ws.Range("A3:C3"))
$ws.Range("A2:C2").Copy($ws.Range("A4:C4"))
$ws.Range("A2:C2").Copy($ws.Range("A5:C5"))

# Append row 7 as a duplicate of row 6 (the original row 3's data, now shifted).
$ws.Range("A6:C6").Copy($ws.Range("A7:C7"))
